$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Locate the data row for participant 5401 -- the new participant row
# (5390, 60.4%, 39.6%) belongs immediately before it, keeping the
# table's ascending participant-id order (... 5328, 5390, 5401 ...).
$afterRow = $null
for ($i = 1; $i -le $t.Rows.Count; $i++) {
  $firstCellText = $t.Rows.Item($i).Cells.Item(1).Range.Text
  if ($firstCellText -like "5401*") {
    $afterRow = $t.Rows.Item($i)
    break
  }
}

$newRow = $t.Rows.Add($afterRow)
$newRow.Cells.Item(1).Range.Text = "5390"
$newRow.Cells.Item(2).Range.Text = "60.4"
$newRow.Cells.Item(3).Range.Text = "39.6"
